$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("D9").Value = "데이터 사이언스랑 경제학이랑 무슨 관련이 있나요?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/ds-econ-how-close/#utm_source=rss&utm_medium=rss&utm_campaign=ds-econ-how-close"

# Row 26
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 32
$ws.Range("D32").Value = "tensorflow_decision_forests 를 이용해서 손쉽게 RandomForest, GBM 사용하기"
$ws.Range("E32").Value = "https://dodonam.tistory.com/377"

# Row 37
$ws.Range("D37").Value = "[Paper Review] Self-Adaptive Forecasting for Improved Deep Learning on Non-Stationary Time-Series"

# Row 51
$ws.Range("D51").Value = "[sqlite3] csv 파일을 sqlite로 가져오는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/1328"
